$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "sv"
$ws.Range("J2").Value = "Statement-opinion"

$ws.Range("I4").Value = "sd"
$ws.Range("J4").Value = "Statement-non-opinion"

$ws.Range("I6").Value = "sv"
$ws.Range("J6").Value = "Statement-opinion"

$ws.Range("I9").Value = "aa"
$ws.Range("J9").Value = "Agree/Accept"

$ws.Range("I12").Value = "%"
$ws.Range("J12").Value = "Uninterpretable"

$ws.Range("I29").Value = "sv"
$ws.Range("J29").Value = "Statement-opinion"
